$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------
# survey sheet: mark the grid rows as hidden-in-contents (F column = TRUE)
# ---------------------------------------------------------------------
$rows = 60,64,68,72,76,80,84,88
foreach ($r in $rows) {
    $survey.Cells.Item($r, 6).Value = $true
}

# widen column E on the survey sheet
$survey.Columns.Item(5).ColumnWidth = 50

# ---------------------------------------------------------------------
# choices sheet: data_value / display.text columns become the string "yes"
# ---------------------------------------------------------------------
foreach ($r in 2..7) {
    $choices.Cells.Item($r, 2).Value = "yes"
    $choices.Cells.Item($r, 3).Value = "yes"
}

# ---------------------------------------------------------------------
# window / sheet view state: "choices" tab becomes the active / selected tab,
# "survey" loses the tabSelected flag and scrolls/selects differently
# ---------------------------------------------------------------------
$survey.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1
$survey.Range("F40").Select()

$choices.Activate()
$choices.Range("C8").Select()

Write-Output "done"
